# "second draft / remove fig6 for now"
#
# Resizes/repositions several caption / label textboxes on slide 1 (they all
# shrink in width and several shift upward, since a figure was removed), and
# merges a run-split piece of text back into a single run.
#
# NOTE on numeric precision: PowerPoint's Shape.Left/.Top/.Width/.Height are
# exposed as single-precision (float32) point values, so not every EMU value
# is exactly reachable from an arbitrary point value. The literals below were
# chosen (by searching nearby float32 values) to reproduce the exact target
# EMU after PowerPoint's internal pt*12700 conversion, landing exactly on
# target in all but a few cases where the target EMU sits strictly between
# two adjacent representable float32 points (those land within 1 EMU, i.e.
# far less than 1/12700 of a point).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($sh in $s.Shapes) {
    switch ($sh.Id) {
        18 {
            # "B. Articulation (sound duration)" - narrower only
            $sh.Width = 732.8467407226562
        }
        19 {
            # "C. Dynamics (velocity profiles)" - narrower + slightly shorter
            $sh.Width  = 740.0204467773438
            $sh.Height = 79.50402069091797
        }
        20 {
            # "(I) Extract performances features from dataset" was split
            # across two runs with identical formatting; merge them into one
            # run without disturbing the shape's autofit-computed height.
            $tr = $sh.TextFrame.TextRange
            $r1 = $tr.Runs(1, 1)
            $r2 = $tr.Runs(2, 1)
            $r1.Text = "(I) Extract performances features from dataset"
            $r2.Text = ""
        }
        30 {
            # "D. Articulation (sound duration)" - shifts down 1 EMU, narrower
            $sh.Top    = 1023.3565673828125
            $sh.Width  = 740.0204467773438
            $sh.Height = 79.50394439697266
        }
        32 {
            # "E. Dynamics (velocity profiles)" - shifts down 1 EMU, narrower
            $sh.Top    = 1525.5059814453125
            $sh.Width  = 732.8467407226562
            $sh.Height = 79.973388671875
        }
        33 {
            # "None recordings (16 instances) - Combine A, B, C..." moves up
            $sh.Top = 2145.698974609375
        }
        47 {
            # "16 instances for D were generated..." - narrower only
            $sh.Width = 732.8467407226562
        }
        8 {
            # "(II) Combine extracted features" moves up
            $sh.Top = 1905.337890625
        }
        9 {
            # "None recordings (4 instances) - none_1, none_2, ..." moves up
            $sh.Top = 2145.534423828125
        }
        10 {
            # "(III) Selected recordings ... (Fig. 7 - Fig. 10 for details)" moves up
            $sh.Top = 1903.62353515625
        }
    }
}
